# Add Files via upload - timesheet updates for 이미정 (sheet4) and 탁재인 (sheet6)

$wb = $excel.ActiveWorkbook

# --- 이미정 sheet: fill in new time-log rows 9-12 (previously blank template rows) ---
$wsLee = $wb.Worksheets.Item("이미정")

# Row 9 : 9월 17일, 21:00 - 23:00, interruption 0, delta 120 min, 회의
$wsLee.Range("A9").Value = "9월 17일"
$wsLee.Range("B9").Value = 0.875
$wsLee.Range("C9").Value = 0.95833333333333337
$wsLee.Range("D9").Value = 0
$wsLee.Range("E9").Value = 120
$wsLee.Range("F9").Value = "회의"

# Row 10 : 9월 19일, 22:00 - 23:00, interruption 0, delta 60 min, 회의
$wsLee.Range("A10").Value = "9월 19일"
$wsLee.Range("B10").Value = 0.91666666666666663
$wsLee.Range("C10").Value = 0.95833333333333337
$wsLee.Range("D10").Value = 0
$wsLee.Range("E10").Value = 60
$wsLee.Range("F10").Value = "회의"

# Row 11 : 9월 20일, 22:30 - 01:30, interruption 0, delta 180 min, 회의
$wsLee.Range("A11").Value = "9월 20일"
$wsLee.Range("B11").Value = 0.9375
$wsLee.Range("C11").Value = 0.0625
$wsLee.Range("D11").Value = 0
$wsLee.Range("E11").Value = 180
$wsLee.Range("F11").Value = "회의"

# Row 12 : 9월 22일, 22:00 - 12:00, interruption 0, delta 120 min, 회의
$wsLee.Range("A12").Value = "9월 22일"
$wsLee.Range("B12").Value = 0.91666666666666663
$wsLee.Range("C12").Value = 0.5
$wsLee.Range("D12").Value = 0
$wsLee.Range("E12").Value = 120
$wsLee.Range("F12").Value = "회의"

# --- 탁재인 sheet: correct the time entry on row 10 (9월 20일) ---
$wsTak = $wb.Worksheets.Item("탁재인")
$wsTak.Range("B10").Value = 0.91666666666666663
$wsTak.Range("C10").Value = 0.020833333333333332
$wsTak.Range("E10").Value = 150

# --- TOTAL sheet: page setup now carries an explicit paper size / orientation ---
$wsTotal = $wb.Worksheets.Item("TOTAL")
$wsTotal.PageSetup.PaperSize = 9
$wsTotal.PageSetup.Orientation = 1

# --- restore the cursor/selection state left behind in the saved view ---
# (order matters: the sheet selected last becomes the active tab, and 탁재인
# is the sheet that should stay active/selected, matching the saved workbook)
$wsKim = $wb.Worksheets.Item("김혜민")
$wsKim.Range("E14").Select()

$wsLee.Range("E15").Select()

$wsTak.Range("D16").Select()
